# "su modif Data para regresiones en Preprod"
# Update the QA data row (row 11, "Baioni Alejandro Luis" record) on Hoja1
# so it points at "Agustin Seisdedos" instead, with refreshed
# Documento/PAS/Answer numbers, then leave I11:J11 selected like the
# author did after editing those cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Documento (G11): 24741860 -> 24741861
$ws.Cells.Item(11, 7).Value = 24741861

# PAS (I11): 234 -> 6254
$ws.Cells.Item(11, 9).Value = 6254

# Answer (J11): "Baioni Alejandro Luis" -> "Agustin Seisdedos"
$ws.Cells.Item(11, 10).Value = "Agustin Seisdedos"

# NumeroCalle (N11): 305 -> 306
$ws.Cells.Item(11, 14).Value = 306

# Match the author's final selection after editing the row.
$ws.Range("I11:J11").Select()
